# Insert a new weekly record at row 27, shifting all existing rows
# (old rows 27-64) down by one, to become rows 28-65.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 27 and below down by one row.
$ws.Rows("27:27").Insert()

# Populate the newly inserted row 27 with the new record's data.
$ws.Cells.Item(27, 1).Value = 3
$ws.Cells.Item(27, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(27, 3).Value = "Coquimbo"
$ws.Cells.Item(27, 4).Value = 44771
$ws.Cells.Item(27, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(27, 5).Value = 5
$ws.Cells.Item(27, 6).Value = 100112035
$ws.Cells.Item(27, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 55
$ws.Cells.Item(27, 11).Value = 15000
$ws.Cells.Item(27, 12).Value = 15000
$ws.Cells.Item(27, 13).Value = 15000
$ws.Cells.Item(27, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(27, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(27, 16).Value = 1000
$ws.Cells.Item(27, 17).Value = 15
$ws.Cells.Item(27, 18).Value = "Hortaliza"
